$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.03659999999999999
$ws.Range("E2").Value = 0.04295
$ws.Range("G2").Value = 0.1468362647840464
$ws.Range("H2").Value = 0.1468362647840464
$ws.Range("I2").Value = 0.1655813198628609
$ws.Range("J2").Value = 0.1264723454239262
$ws.Range("K2").Value = 475.9
$ws.Range("L2").Value = 0.09654514839835272
$ws.Range("M2").Value = 341.242
$ws.Range("N2").Value = 0.02997663305105591
$ws.Range("O2").Value = 0.7170455978146669
$ws.Range("P2").Value = 313.742
$ws.Range("Q2").Value = 0.02756087705119646
$ws.Range("R2").Value = 0.6592603488127758
$ws.Range("S2").Value = 27.49999999999999
$ws.Range("T2").Value = 0.08058796982786406
$ws.Range("U2").Value = 327
$ws.Range("V2").Value = 0.02872553498014688
$ws.Range("W2").Value = 0.1785441822737376
$ws.Range("X2").Value = 0.04308125810136546
$ws.Range("Y2").Value = 0.1354629241723722
$ws.Range("Z2").Value = 1.589584005159626
$ws.Range("AA2").Value = 0.2207215824213049
$ws.Range("AB2").Value = 0.04008440751776571
$ws.Range("AC2").Value = 0.1806371749035393
$ws.Range("AD2").Value = 962
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 962
$ws.Range("AG2").Value = 635
$ws.Range("AH2").Value = 0.07792249870399172
$ws.Range("AI2").Value = 0.2542955326460481
$ws.Range("AJ2").Value = 0.05283477276887492
$ws.Range("AK2").Value = 0.1837384259259259
$ws.Range("AL2").Value = 32.78
$ws.Range("AM2").Value = 32.78
$ws.Range("AN2").Value = 1.074260189838079
$ws.Range("AO2").Value = 24.8993288590604
$ws.Range("AP2").Value = 0.7091010608598548
$ws.Range("AQ2").Value = 24.8993288590604
# Row 3
$ws.Range("B3").Value = "Alm. Brand A/S (CPSE:ALMB)"
$ws.Range("D3").Value = 0.0189
$ws.Range("E3").Value = 0.0295
$ws.Range("G3").Value = 0.08637023718968261
$ws.Range("H3").Value = 0.08637023718968261
$ws.Range("I3").Value = 0.1890602309660466
$ws.Range("J3").Value = 0.1471726348837488
$ws.Range("K3").Value = 91
$ws.Range("L3").Value = 0.06292787497406818
$ws.Range("M3").Value = 86.142
$ws.Range("N3").Value = 0.04638022936520756
$ws.Range("O3").Value = 0.9466153846153845
$ws.Range("P3").Value = 72.842
$ws.Range("Q3").Value = 0.03921929682872988
$ws.Range("R3").Value = 0.8004615384615384
$ws.Range("S3").Value = 13.3
$ws.Range("T3").Value = 0.1543962294815537
$ws.Range("U3").Value = 222.7
$ws.Range("V3").Value = 0.1199052387874872
$ws.Range("W3").Value = 0.1355376824545725
$ws.Range("X3").Value = 0.04433154243519008
$ws.Range("Y3").Value = 0.09120614001938246
$ws.Range("Z3").Value = 1.79595131644312
$ws.Range("AA3").Value = 0.2643148873638714
$ws.Range("AB3").Value = 0.04003784737526019
$ws.Range("AC3").Value = 0.2242770399886112
$ws.Range("AD3").Value = 318.2
$ws.Range("AF3").Value = 318.2
$ws.Range("AG3").Value = 95.5
$ws.Range("AH3").Value = 0.1462652263847391
$ws.Range("AI3").Value = 0.283423888839405
$ws.Range("AJ3").Value = 0.04890413764850472
$ws.Range("AK3").Value = 0.1061111111111111
$ws.Range("AL3").Value = 8.68
$ws.Range("AM3").Value = 8.68
$ws.Range("AN3").Value = 1.013053167780961
$ws.Range("AO3").Value = 31.49769585253456
$ws.Range("AP3").Value = 0.3040432983126393
$ws.Range("AQ3").Value = 31.49769585253456
# Row 4
$ws.Range("B4").Value = "Tryg A/S (CPSE:TRYG)"
$ws.Range("D4").Value = 0.05429999999999999
$ws.Range("E4").Value = 0.0564
$ws.Range("G4").Value = 0.1719395957740009
$ws.Range("H4").Value = 0.1719395957740009
$ws.Range("I4").Value = 0.1558337161231052
$ws.Range("J4").Value = 0.1167464111838791
$ws.Range("K4").Value = 384.9
$ws.Range("L4").Value = 0.1105018373909049
$ws.Range("M4").Value = 255.1
$ws.Range("N4").Value = 0.02677849742292391
$ws.Range("O4").Value = 0.6627695505326059
$ws.Range("P4").Value = 240.9
$ws.Range("Q4").Value = 0.02528788721749263
$ws.Range("R4").Value = 0.6258768511301638
$ws.Range("S4").Value = 14.19999999999999
$ws.Range("T4").Value = 0.05566444531556248
$ws.Range("U4").Value = 104.3
$ws.Range("V4").Value = 0.01094863693144243
$ws.Range("W4").Value = 0.2215506820929028
$ws.Range("X4").Value = 0.04183097376754084
$ws.Range("Y4").Value = 0.1797197083253619
$ws.Range("Z4").Value = 1.517205331474867
$ws.Range("AA4").Value = 0.1771282774787385
$ws.Range("AB4").Value = 0.04013096766027122
$ws.Range("AC4").Value = 0.1369973098184673
$ws.Range("AD4").Value = 643.8
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 643.8
$ws.Range("AG4").Value = 539.5
$ws.Range("AH4").Value = 0.06330321235779393
$ws.Range("AI4").Value = 0.2420027816411683
$ws.Range("AJ4").Value = 0.05359732957142006
$ws.Range("AK4").Value = 0.2110719874804382
$ws.Range("AL4").Value = 24.1
$ws.Range("AM4").Value = 24.1
$ws.Range("AN4").Value = 1.107327141382869
$ws.Range("AO4").Value = 22.52282157676348
$ws.Range("AP4").Value = 0.9279325765393878
$ws.Range("AQ4").Value = 22.52282157676348
